# Tchoukball Scoreboard template update
# Commit: "80 - add option to switch between millisecond and seconds timer"
# The accompanying data change replaces the demo fixture list (old
# Malaysian-state / age-group sample data) with a new set of matches,
# and adds an 8th match row. Each row's Category gets a colour-coded
# CategoryColor cell (column G) so the scoreboard app can pick up the
# right highlight colour for the category.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colour per category (OLE/COM color = R + G*256 + B*65536)
#   U18 BOY      -> 00B0F0 (blue)
#   OPEN FEMALE  -> 92D050 (green)
#   OPEN MALE    -> FFFF00 (yellow)
#   U18 GIRL     -> FF9900 (orange)
$colorMap = @{
    "U18 BOY"     = 15773696
    "OPEN FEMALE" = 5296274
    "OPEN MALE"   = 65535
    "U18 GIRL"    = 39423
}

$rows = @(
    @{ Row=2; Match="A01"; Home="SELANGOR"; Away="JOHOR";    Category="U18 BOY" },
    @{ Row=3; Match="A02"; Home="PENANG";   Away="JOHOR 1";  Category="OPEN FEMALE" },
    @{ Row=4; Match="A03"; Home="JOHOR 2";  Away="SARAWAK";  Category="OPEN MALE" },
    @{ Row=5; Match="A04"; Home="JOHOR";    Away="SARAWAK";  Category="U18 BOY" },
    @{ Row=6; Match="A05"; Home="JOHOR 1";  Away="SELANGOR"; Category="OPEN MALE" },
    @{ Row=7; Match="A06"; Home="JOHOR 1";  Away="JOHOR 2";  Category="U18 GIRL" },
    @{ Row=8; Match="A07"; Home="SELANGOR"; Away="JOHOR 2";  Category="OPEN MALE" },
    @{ Row=9; Match="A08"; Home="SARAWAK";  Away="JOHOR 1";  Category="OPEN MALE" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.Match
    $ws.Range("B$n").Value = $r.Home
    $ws.Range("C$n").Value = $r.Away
    $ws.Range("D$n").Value = "12:00"
    $ws.Range("E$n").Value = "02:00"
    $ws.Range("F$n").Value = $r.Category
    $ws.Range("G$n").Interior.Color = $colorMap[$r.Category]
}

# Restore the cursor to where the author left it after editing.
$ws.Range("F12").Select()
